$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D4").Value = 49055
$ws.Range("E4").Value = 2852
$ws.Range("D5").Value = 50112
$ws.Range("E5").Value = 4026
$ws.Range("D6").Value = 51175
$ws.Range("E6").Value = 5300
$ws.Range("D7").Value = 52234
$ws.Range("E7").Value = 6478
$ws.Range("D8").Value = 53293
$ws.Range("E8").Value = 7587
$ws.Range("C9").Value = 45741
$ws.Range("D9").Value = 54325
$ws.Range("E9").Value = 8584
$ws.Range("D10").Value = 55330
$ws.Range("E10").Value = 9479
$ws.Range("D11").Value = 56277
$ws.Range("E11").Value = 10249
$ws.Range("D12").Value = 57146
$ws.Range("E12").Value = 10909
$ws.Range("D13").Value = 57939
$ws.Range("E13").Value = 11480
$ws.Range("D14").Value = 58633
$ws.Range("E14").Value = 11952
$ws.Range("D15").Value = 59232
$ws.Range("E15").Value = 12322
$ws.Range("D16").Value = 59752
$ws.Range("E16").Value = 12577
$ws.Range("D17").Value = 60185
$ws.Range("E17").Value = 12715
$ws.Range("D18").Value = 60574
$ws.Range("E18").Value = 12786
$ws.Range("C19").Value = 48121
$ws.Range("D19").Value = 60940
$ws.Range("E19").Value = 12820
$ws.Range("D20").Value = 61317
$ws.Range("E20").Value = 12886
$ws.Range("D21").Value = 61723
$ws.Range("E21").Value = 13041
$ws.Range("D22").Value = 62181
$ws.Range("E22").Value = 13312
$ws.Range("C24").Value = 52369
$ws.Range("D24").Value = 52773
$ws.Range("E24").Value = 403
$ws.Range("D25").Value = 53233
$ws.Range("E25").Value = 803
$ws.Range("D26").Value = 53749
$ws.Range("E26").Value = 1420
$ws.Range("D27").Value = 54321
$ws.Range("E27").Value = 2194
$ws.Range("C28").Value = 51766
$ws.Range("D28").Value = 54931
$ws.Range("E28").Value = 3165
$ws.Range("D29").Value = 55547
$ws.Range("E29").Value = 4226
$ws.Range("C30").Value = 50917
$ws.Range("D30").Value = 56193
$ws.Range("E30").Value = 5276
$ws.Range("D31").Value = 56819
$ws.Range("E31").Value = 6258
$ws.Range("D32").Value = 57406
$ws.Range("E32").Value = 7207
$ws.Range("D33").Value = 57928
$ws.Range("E33").Value = 7990
$ws.Range("C34").Value = 49716
$ws.Range("D34").Value = 58391
$ws.Range("E34").Value = 8675
$ws.Range("D35").Value = 58778
$ws.Range("E35").Value = 9234
$ws.Range("C36").Value = 49424
$ws.Range("D36").Value = 59083
$ws.Range("E36").Value = 9659
$ws.Range("D37").Value = 59331
$ws.Range("E37").Value = 9891
$ws.Range("D38").Value = 59522
$ws.Range("E38").Value = 9975
$ws.Range("D39").Value = 59689
$ws.Range("E39").Value = 9978
$ws.Range("D40").Value = 59849
$ws.Range("E40").Value = 9930
$ws.Range("D41").Value = 60004
$ws.Range("E41").Value = 9835
$ws.Range("D42").Value = 60182
$ws.Range("E42").Value = 9717
$ws.Range("C43").Value = 50777
$ws.Range("D43").Value = 60374
$ws.Range("E43").Value = 9597
$ws.Range("C44").Value = 16556
$ws.Range("E44").Value = 1
$ws.Range("C45").Value = 16581
$ws.Range("E45").Value = 998
$ws.Range("C46").Value = 16636
$ws.Range("D46").Value = 18647
$ws.Range("E46").Value = 2011
$ws.Range("C47").Value = 16655
$ws.Range("D47").Value = 19746
$ws.Range("E47").Value = 3091
$ws.Range("C48").Value = 16666
$ws.Range("D48").Value = 20847
$ws.Range("E48").Value = 4181
$ws.Range("D49").Value = 21949
$ws.Range("E49").Value = 5285
$ws.Range("D50").Value = 23031
$ws.Range("E50").Value = 6381
$ws.Range("C51").Value = 16637
$ws.Range("D51").Value = 24103
$ws.Range("E51").Value = 7466
$ws.Range("C52").Value = 16616
$ws.Range("D52").Value = 25145
$ws.Range("E52").Value = 8529
$ws.Range("C53").Value = 16590
$ws.Range("D53").Value = 26173
$ws.Range("E53").Value = 9583
$ws.Range("D54").Value = 27167
$ws.Range("E54").Value = 10594
$ws.Range("C55").Value = 16565
$ws.Range("D55").Value = 28127
$ws.Range("E55").Value = 11561
$ws.Range("C56").Value = 16565
$ws.Range("D56").Value = 29061
$ws.Range("E56").Value = 12496
$ws.Range("C57").Value = 16575
$ws.Range("D57").Value = 29976
$ws.Range("E57").Value = 13401
$ws.Range("C58").Value = 16591
$ws.Range("D58").Value = 30886
$ws.Range("E58").Value = 14295
$ws.Range("D59").Value = 31782
$ws.Range("E59").Value = 15169
$ws.Range("C60").Value = 16643
$ws.Range("D60").Value = 32658
$ws.Range("E60").Value = 16015
$ws.Range("C61").Value = 16689
$ws.Range("D61").Value = 33531
$ws.Range("E61").Value = 16842
$ws.Range("C62").Value = 16744
$ws.Range("D62").Value = 34404
$ws.Range("E62").Value = 17660
$ws.Range("D63").Value = 35275
$ws.Range("E63").Value = 18471
$ws.Range("C64").Value = 16866
$ws.Range("D64").Value = 36130
$ws.Range("E64").Value = 19264
$ws.Range("D67").Value = 30890
$ws.Range("D69").Value = 31130
$ws.Range("E70").Value = -20
$ws.Range("E71").Value = 50
$ws.Range("D73").Value = 31594
$ws.Range("E73").Value = 381
$ws.Range("D74").Value = 31672
$ws.Range("E74").Value = 637
$ws.Range("D75").Value = 31717
$ws.Range("E75").Value = 865
$ws.Range("D76").Value = 31726
$ws.Range("E76").Value = 1000
$ws.Range("D77").Value = 31699
$ws.Range("E77").Value = 1058
$ws.Range("D78").Value = 31633
$ws.Range("E78").Value = 1009
$ws.Range("D79").Value = 31553
$ws.Range("E79").Value = 994
$ws.Range("E80").Value = 968
$ws.Range("D81").Value = 31322
$ws.Range("E81").Value = 914
$ws.Range("D82").Value = 31190
$ws.Range("E82").Value = 740
$ws.Range("D83").Value = 31069
$ws.Range("E85").Value = 13
$ws.Range("C86").Value = 13249
$ws.Range("E86").Value = -1
$ws.Range("C87").Value = 13275
$ws.Range("E87").Value = -50
$ws.Range("C88").Value = 13336
$ws.Range("E88").Value = -146
$ws.Range("C89").Value = 13396
$ws.Range("E89").Value = -256
$ws.Range("C90").Value = 13442
$ws.Range("E90").Value = -351
$ws.Range("C91").Value = 13513
$ws.Range("E91").Value = -489
$ws.Range("C92").Value = 13578
$ws.Range("E92").Value = -634
$ws.Range("C93").Value = 13612
$ws.Range("D93").Value = 12866
$ws.Range("E93").Value = -746
$ws.Range("C94").Value = 13570
$ws.Range("E94").Value = -807
$ws.Range("C95").Value = 13502
$ws.Range("E95").Value = -838
$ws.Range("C96").Value = 13420
$ws.Range("E96").Value = -864
$ws.Range("C97").Value = 13362
$ws.Range("E97").Value = -927
$ws.Range("C98").Value = 13318
$ws.Range("D98").Value = 12321
$ws.Range("E98").Value = -997
$ws.Range("C99").Value = 13301
$ws.Range("E99").Value = -1110
$ws.Range("C100").Value = 13235
$ws.Range("E100").Value = -1155
$ws.Range("C101").Value = 13143
$ws.Range("E101").Value = -1182
$ws.Range("C102").Value = 13061
$ws.Range("E102").Value = -1209
$ws.Range("C103").Value = 13033
$ws.Range("E103").Value = -1284
$ws.Range("C104").Value = 13033
$ws.Range("D104").Value = 11635
$ws.Range("E104").Value = -1399
$ws.Range("C105").Value = 13037
$ws.Range("D105").Value = 11542
$ws.Range("E105").Value = -1496
$ws.Range("C106").Value = 13053
$ws.Range("E106").Value = -1610
